# Case_7_44 diagnostic.xlsx: write the B1/A2 numeric flags, the
# B2 "disconnected_elements" label, and the shared bold/bordered/
# centered style used by B1 and A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the style once on B1 ...
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160

# ... then copy its format onto A2, so both cells share a single
# consolidated cell style (rather than each accumulating its own
# chain of intermediate style records).
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
